# Apply the commit "changed pwd in config":
#  - SQL sheet: insert a new Q&A row (21: new srch_consol_tbl query) above the
#    last (email-lookup) row, which pushes that row from 22 -> 23.
#  - Login sheet: update the OptumPwd test value (F4) to a new password.
#  - Switch the active/selected tab from SQL to Login, and update the
#    remembered selection on each sheet.

$wb  = $excel.ActiveWorkbook
$sql   = $wb.Worksheets.Item("SQL")
$login = $wb.Worksheets.Item("Login")

# --- SQL sheet: insert new row 22 (shifts old row 22 -> row 23) -------------
$sql.Rows(22).Insert()

$sql.Range("A22").Value = "'21"

$query = "select * from ole.srch_consol_tbl" + "`n" + `
         "where PROV_TAX_ID_NBR ='{`$tin}'" + "`n" + `
         "and CP_PAY_METH_CD <> 'CHK'" + "`n" + `
         "and CP_SETL_DT>='2018-03-30' and CP_SETL_DT<='2018-05-29'" + "`n" + `
         "order by CP_SETL_DT desc"
$sql.Range("B22").Value = $query

# Match the authored row height (5 wrapped lines @ 14.4pt).
$sql.Rows(22).RowHeight = 72

# --- Login sheet: change the stored OptumPwd test value ---------------------
$login.Range("F4").Value = "Test_1234"

# --- View state: Login becomes the active tab; remember each sheet's
#     selection (SQL -> B5, Login keeps F4) and drop SQL's scrolled
#     topLeftCell. -----------------------------------------------------------
$sql.Activate()
$sql.Range("B5").Select()
$login.Activate()
$login.Range("F4").Select()
